$wb = $excel.ActiveWorkbook

# Rename the first sheet
$ws = $wb.Worksheets.Item(1)
$ws.Name = "ValidLogin"

# Set header row
$ws.Range("A1").Value = "UserName"
$ws.Range("B1").Value = "Password"

# Set data row
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "manager"

# Update selection to match target
$ws.Range("G16").Select()
